$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column widths for D, E, F
$ws.Columns.Item(4).ColumnWidth = 23.59
$ws.Columns.Item(5).ColumnWidth = 18.74
$ws.Columns.Item(6).ColumnWidth = 30.65

# Row 1: D1,E1,F1 = 5
$ws.Range("D1").Value = 5
$ws.Range("E1").Value = 5
$ws.Range("F1").Value = 5

# Row 2: C2 stays -5 (already present); D2,E2,F2 = -5
$ws.Range("D2").Value = -5
$ws.Range("E2").Value = -5
$ws.Range("F2").Value = -5

# Row 3: D3,E3,F3 = 5.1234567
$ws.Range("D3").Value = 5.1234567
$ws.Range("E3").Value = 5.1234567
$ws.Range("F3").Value = 5.1234567

# Row 4: D4,E4,F4 = 5.12345678987654
$ws.Range("D4").Value = 5.12345678987654
$ws.Range("E4").Value = 5.12345678987654
$ws.Range("F4").Value = 5.12345678987654

# Row 5: D5,E5,F5 = 12345678987654300
$ws.Range("D5").Value = 12345678987654300
$ws.Range("E5").Value = 12345678987654300
$ws.Range("F5").Value = 12345678987654300

# Selection moves to F4
$ws.Range("F4").Select()
